$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (values 9 / 15, bordered/centered style) is being
# removed, and every other column shifts one place to the left
# (old B->A, C->B, D->C, E->D, F->E). The values that used to live in
# column F (9 / 15) end up in the new column E, but without the old
# column A's style.

# Shift everything left by deleting column A entirely.
$ws.Range("A1").EntireColumn.Delete()

# After the delete, old F (now E) already holds 9 / 15 from the shift,
# so nothing else to move - but make sure new E2/E3 carry plain (no
# special style) formatting, matching the target.
$ws.Range("E2:E3").Style = "Normal"
